$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.601.18'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.739.67'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.23'
$ws.Range("E5").Value = '  +0.67%  '


$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4906'
$ws.Range("E7").Value = '  +2.49%  '

$ws.Range("E8").Value = '  -0.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06281'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.747.61'
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07049'
$ws.Range("E11").Value = '  -1.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.71'
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6135'
$ws.Range("E13").Value = '  -0.48%  '

$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.00'
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.608.75'
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007289'
$ws.Range("E19").Value = '  +5.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.973.32'
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.567'
$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.712'
$ws.Range("E23").Value = '  -2.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.275'
$ws.Range("E24").Value = '  -0.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.92'
$ws.Range("E25").Value = '  +1.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.42'
$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.421'
$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("E28").Value = '  -2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.39'
$ws.Range("E29").Value = '  +0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.027'
$ws.Range("E30").Value = '  +1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08058'
$ws.Range("E31").Value = '  +0.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.724'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04606'
$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.013'
$ws.Range("E35").Value = '  +2.55%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6400'
$ws.Range("E36").Value = '  +0.62%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.070'
$ws.Range("E37").Value = '  +1.00%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9050'
$ws.Range("E38").Value = '  -3.06%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.430'
$ws.Range("E39").Value = '  +0.80%  '

$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.003'
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01504'
$ws.Range("E41").Value = '  +0.40%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.18'
$ws.Range("E42").Value = '  -4.39%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.430'
$ws.Range("E43").Value = '  -3.59%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3914'
$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.866'
$ws.Range("E45").Value = '  -1.88%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1185'
$ws.Range("E46").Value = '  -0.51%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05397'
$ws.Range("E47").Value = '  +1.40%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.61'
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.807'
$ws.Range("E49").Value = '  -0.59%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.258'
$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.78'
$ws.Range("E51").Value = '  +0.64%  '
